$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format before writing, so values like
# "554.74" or "1.00" are stored as literal strings (matching the workbook's
# original inlineStr cells) rather than being auto-coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.833.15"
$ws.Range("D3").Value = "2.425.70"
$ws.Range("D5").Value = "552.44"
$ws.Range("D6").Value = "160.30"
$ws.Range("D8").Value = "0.495"
$ws.Range("D9").Value = "2.425.67"
$ws.Range("D10").Value = "0.146"
$ws.Range("D11").Value = "0.163"
$ws.Range("D12").Value = "0.331"
$ws.Range("D13").Value = "4.70"
$ws.Range("D14").Value = "2.881.42"
$ws.Range("D15").Value = "67.827.85"
$ws.Range("D16").Value = "0.0000165"
$ws.Range("D17").Value = "22.81"
$ws.Range("D18").Value = "2.435.60"
$ws.Range("D19").Value = "10.67"
$ws.Range("D20").Value = "335.84"
$ws.Range("D21").Value = "6.98"
$ws.Range("D22").Value = "3.67"
$ws.Range("D23").Value = "1.00"
$ws.Range("D25").Value = "65.68"
$ws.Range("D27").Value = "3.57"
$ws.Range("D28").Value = "1.00"
$ws.Range("D29").Value = "7.92"
$ws.Range("D30").Value = "0.0₃0811"
$ws.Range("D31").Value = "7.00"
$ws.Range("D33").Value = "415.68"
$ws.Range("D34").Value = "1.12"
$ws.Range("D35").Value = "1.60"
$ws.Range("D36").Value = "157.53"
$ws.Range("D39").Value = "0.107"
$ws.Range("D40").Value = "17.56"
$ws.Range("D41").Value = "0.297"
$ws.Range("D42").Value = "4.26"
$ws.Range("D43").Value = "1.44"
$ws.Range("D44").Value = "1.06"
$ws.Range("D45").Value = "2.01"
$ws.Range("D46").Value = "132.01"
$ws.Range("D47").Value = "3.26"
$ws.Range("D48").Value = "0.0710"
$ws.Range("D49").Value = "0.468"
$ws.Range("D50").Value = "0.551"
$ws.Range("D51").Value = "0.0898"

# Restore the default "Normal" style on the Price column so no stray
# number-format style lingers on cells (keeps styles.xml unchanged).
$ws.Range("D2:D51").Style = "Normal"

# Volume(1h) percentages: the surrounding spaces keep these safely text,
# so they can be written directly.
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  -2.33%  "
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("E10").Value = "  -7.35%  "
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("E12").Value = "  -6.12%  "
$ws.Range("E13").Value = "  -4.00%  "
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("E16").Value = "  -5.36%  "
$ws.Range("E17").Value = "  -5.72%  "
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("E19").Value = "  -4.04%  "
$ws.Range("E20").Value = "  -2.77%  "
$ws.Range("E21").Value = "  -5.23%  "
$ws.Range("E22").Value = "  -4.71%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -5.23%  "
$ws.Range("E25").Value = "  -5.30%  "
$ws.Range("E27").Value = "  -7.65%  "
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  -8.44%  "
$ws.Range("E30").Value = "  -6.80%  "
$ws.Range("E31").Value = "  -8.89%  "
$ws.Range("E33").Value = "  -5.99%  "
$ws.Range("E34").Value = "  -5.88%  "
$ws.Range("E35").Value = "  -5.58%  "
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  -5.20%  "
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("E41").Value = "  -4.98%  "
$ws.Range("E42").Value = "  -6.84%  "
$ws.Range("E43").Value = "  -7.99%  "
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("E45").Value = "  -6.64%  "
$ws.Range("E46").Value = "  -5.40%  "
$ws.Range("E47").Value = "  -4.84%  "
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("E49").Value = "  -8.48%  "
$ws.Range("E50").Value = "  -3.60%  "
$ws.Range("E51").Value = "  -2.28%  "
